# Add files via upload
#
# The sheet gains a new data row (row 3) with a date/number/number/text
# record, and the existing date column's number format gains an explicit
# time component (the "Data" column now formats as date + time instead of
# date only), while the newly added date cell keeps the original
# date-only display.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing date cell (A2) switches to a date+time format.
# Do this before formatting the new row so the style slots line up the
# same way Excel would allocate them.
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row of data.
$ws.Range("A3").Value = 45619
$ws.Range("A3").NumberFormat = "YYYY-MM-DD"
$ws.Range("B3").Value = 63
$ws.Range("C3").Value = 223
$ws.Range("D3").Value = "Mudou-se"
